# "unify the docstrings and the documentation to make sense together"
#
# 1. The "Equipment" sheet's K1 header said "Calibration Period [years]" but
#    the rest of the project (docs / README) calls it a "Calibration Cycle".
#    Re-word the header so the workbook and the documentation agree, while
#    keeping the existing rich-text look (bold label, plain units suffix).
# 2. The workbook used to open on the "Connections" tab with cell H4
#    selected there (and G3 selected back on "Equipment"). Make "Equipment"
#    the active/selected sheet with the default A1 selection on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Equipment")
$ws2 = $wb.Worksheets.Item("Connections")

# --- 1. Reword the calibration header -------------------------------------
$cell = $ws1.Range("K1")
$cell.Value = "Calibration Cycle [years]"

# Restore the two-run rich-text formatting: "Calibration Cycle " (bold)
# followed by "[years]" (regular).
$boldRun = $cell.Characters(1, 18)
$boldRun.Font.Bold = $true
$boldRun.Font.Size = 11
$boldRun.Font.Name = "Calibri"

$restRun = $cell.Characters(19, 7)
$restRun.Font.Bold = $false
$restRun.Font.Size = 11
$restRun.Font.Name = "Calibri"

# --- 2. Fix up the active sheet / selection --------------------------------
# Reset "Connections" selection to A1 and leave it the non-active tab.
$null = $ws2.Activate()
$null = $ws2.Range("A1").Select()

# "Equipment" becomes the active tab, selection reset to A1.
$null = $ws1.Activate()
$null = $ws1.Range("A1").Select()
